$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update credentials
$ws.Range("B1").Value = "Rola.khalaf2@bd.com"
$ws.Range("B2").Value = "Carefusion@5"

# Update test data row (row 5)
$ws.Range("A5").Value = "TestAuto_POC31sep"
$ws.Range("B5").Value = "TestAuto_POC31sep"
$ws.Range("C5").Value = "Facility_POC31sep"
$ws.Range("D5").Value = "Facility_POC31sep"
$ws.Range("E5").Value = "Pharmacy_POC31sep"
$ws.Range("F5").Value = "Pharmacy_POC31sep"
$ws.Range("H5").Value = "Epic1011024"

# Update active selection to H5
$ws.Range("H5").Select()
